$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$shape = $s.Shapes.Item(5)
$tf = $shape.TextFrame
$tr = $tf.TextRange
$para = $tr.Paragraphs(2, 1)
$para.Text = "Il Server TCP è sequenziale, quindi riesce a sostenere molteplici connessioni senza dilatare i tempi di esecuzione."
